$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 40 mirrors the formatting of row 39 (same column styles / wrap / row height).
$ws.Range("A39:H39").Copy() | Out-Null
$ws.Range("A40:H40").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Rows(40).RowHeight = 38

$ws.Range("A40").Value = "KK"
$ws.Range("B40").Value = "TALAK1_BLG_P47"
$ws.Range("C40").Value = "L2100"
$ws.Range("D40").Value = "16-Dec-2025 1:25 PM"
$ws.Range("E40").Value = "FAIL"
$ws.Range("F40").Value = "1. VOLTE setup SR`n2. Check functioning of ViLTE"
$ws.Range("G40").Value = "1. Drive Volte SCMO`n2. Static Video call MO"
$ws.Range("H40").Value = "1. Add a VoLTE short-call drive across all sectors. Do not disconnect the call manually, as it will be counted as a call drop. Each sector must have at least one successful MO session setup.`n2. Do not use WhatsApp for this test. Perform a manual VoLTE video call after running the script."

$ws.Range("E40").Select() | Out-Null
